$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "57.948.27"
$ws.Range("E2").Value = "  -3.99%  "
$ws.Range("D3").Value = "2.290.33"
$ws.Range("E3").Value = "  -5.06%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'535.34"
$ws.Range("E5").Value = "  -4.55%  "
$ws.Range("D6").Value = "'131.68"
$ws.Range("E6").Value = "  -3.02%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -3.29%  "
$ws.Range("D9").Value = "2.289.37"
$ws.Range("E9").Value = "  -4.99%  "
$ws.Range("E10").Value = "  -5.75%  "
$ws.Range("D11").Value = "'5.45"
$ws.Range("E11").Value = "  -3.09%  "
$ws.Range("E12").Value = "  -0.95%  "
$ws.Range("E13").Value = "  -4.96%  "
$ws.Range("E14").Value = "  -4.99%  "
$ws.Range("D15").Value = "2.698.97"
$ws.Range("E15").Value = "  -4.99%  "
$ws.Range("D16").Value = "57.896.80"
$ws.Range("E16").Value = "  -3.91%  "
$ws.Range("E17").Value = "  -4.96%  "
$ws.Range("D18").Value = "2.282.20"
$ws.Range("E18").Value = "  -5.30%  "
$ws.Range("E19").Value = "  -5.92%  "
$ws.Range("E20").Value = "  -6.47%  "
$ws.Range("D21").Value = "'312.84"
$ws.Range("E21").Value = "  -3.84%  "
$ws.Range("E22").Value = "  -6.21%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "'62.96"
$ws.Range("E24").Value = "  -2.54%  "
$ws.Range("E25").Value = "  -5.02%  "
$ws.Range("D26").Value = "'1.01"
$ws.Range("E26").Value = "  +0.85%  "
$ws.Range("D27").Value = "'7.97"
$ws.Range("E27").Value = "  -6.23%  "
$ws.Range("E28").Value = "  -5.29%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'169.89"
$ws.Range("E29").Value = "  -0.57%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'1.71"
$ws.Range("E30").Value = "  -5.98%  "
$ws.Range("D31").Value = "0.0₃0723"
$ws.Range("E31").Value = "  -6.50%  "
$ws.Range("D32").Value = "'1.07"
$ws.Range("E32").Value = "  -1.53%  "
$ws.Range("E33").Value = "  -6.67%  "
$ws.Range("E34").Value = "  -6.26%  "
$ws.Range("D36").Value = "'17.67"
$ws.Range("E36").Value = "  -3.69%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("E38").Value = "  -7.68%  "
$ws.Range("E39").Value = "  -6.85%  "
$ws.Range("D40").Value = "'37.98"
$ws.Range("E40").Value = "  -2.02%  "
$ws.Range("E41").Value = "  -7.25%  "
$ws.Range("E42").Value = "  -3.63%  "
$ws.Range("D43").Value = "'288.39"
$ws.Range("E43").Value = "  -10.84%  "
$ws.Range("E44").Value = "  -4.38%  "
$ws.Range("E45").Value = "  -2.72%  "
$ws.Range("D46").Value = "'0.0498"
$ws.Range("E46").Value = "  -3.73%  "
$ws.Range("E47").Value = "  -3.70%  "
$ws.Range("D48").Value = "'18.22"
$ws.Range("E49").Value = "  -5.29%  "
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("E51").Value = "  -0.58%  "
